$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values to trigger the formula recalculations in column D
$ws.Range("B3").Value = 7.56
$ws.Range("C4").Value = 4.2
$ws.Range("C5").Value = 2.7

# Update the active selection to match the final state
$ws.Range("C6").Select()
